$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.912.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.373.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.472"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.121"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("E11").Value = "  -3.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.950.53"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.380.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("E16").Value = "  -2.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.013.51"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("E18").Value = "  -3.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.92"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.56%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000112"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.190"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.13"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.14%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.04"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "165.81"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.409.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0759"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.48"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -9.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.774"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("E42").Value = "  -2.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.62"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.97%  "

$ws.Range("E44").Value = "  -1.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.441.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.51%  "

$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0259"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("E51").Value = "  -3.51%  "
